$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 576.125
$ws.Range("I4").Value = 177.25
$ws.Range("J4").Value = 975
$ws.Range("K4").Value = 177.25
$ws.Range("L4").Value = 975
$ws.Range("M4").Value = -63.25
$ws.Range("N4").Value = -1203
$ws.Range("H33").Value = 18949.666
$ws.Range("I33").Value = 25180.75
$ws.Range("J33").Value = 6487.5
$ws.Range("K33").Value = 25180.75
$ws.Range("L33").Value = 6487.5
$ws.Range("M33").Value = -24951.75
$ws.Range("N33").Value = -6945.5
$ws.Range("H55").Value = 686.9375
$ws.Range("J55").Value = 913.2857
$ws.Range("L55").Value = 913.2857
$ws.Range("N55").Value = -1341.2857
$ws.Range("H113").Value = 4577.875
$ws.Range("I113").Value = 3584.5
$ws.Range("J113").Value = 4909
$ws.Range("K113").Value = 3584.5
$ws.Range("L113").Value = 4909
$ws.Range("M113").Value = -330.5
$ws.Range("N113").Value = -11417
$ws.Range("H132").Value = 3294.5881
$ws.Range("I132").Value = 2550.8125
$ws.Range("K132").Value = 7652.4375
$ws.Range("M132").Value = -5122.4375
$ws.Range("H135").Value = 1792.0526
$ws.Range("I135").Value = 1203.4667
$ws.Range("J135").Value = 3999.25
$ws.Range("K135").Value = 10831.2003
$ws.Range("L135").Value = 35993.25
$ws.Range("M135").Value = -8296.2003
$ws.Range("N135").Value = -41063.25

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 21402.035
$ws.Range("I32").Value = 22181.23
$ws.Range("J32").Value = 13298.4
$ws.Range("K32").Value = 22181.23
$ws.Range("L32").Value = 13298.4
$ws.Range("M32").Value = -21894.23
$ws.Range("N32").Value = -13872.4
$ws.Range("H61").Value = 2635.6428
$ws.Range("I61").Value = 2635.6428
$ws.Range("K61").Value = 2635.6428
$ws.Range("M61").Value = -2423.6428
$ws.Range("H74").Value = 40707.12
$ws.Range("I74").Value = 40707.12
$ws.Range("K74").Value = 40707.12
$ws.Range("M74").Value = -39833.12
$ws.Range("H77").Value = 40707.12
$ws.Range("I77").Value = 40707.12
$ws.Range("K77").Value = 203535.6
$ws.Range("M77").Value = -199167.6
$ws.Range("H94").Value = 40000
$ws.Range("J94").Value = 40000
$ws.Range("L94").Value = 40000
$ws.Range("N94").Value = -41802
$ws.Range("H96").Value = 0
$ws.Range("J96").Value = 0
$ws.Range("L96").Value = 0
$ws.Range("N96").ClearContents()
$ws.Range("H101").Value = 51332.332
$ws.Range("J101").Value = 51332.332
$ws.Range("L101").Value = 51332.332
$ws.Range("N101").Value = -57822.332
$ws.Range("H106").Value = 0
$ws.Range("J106").Value = 0
$ws.Range("L106").Value = 0
$ws.Range("N106").ClearContents()
$ws.Range("H117").Value = 124999.75
$ws.Range("J117").Value = 124999.75
$ws.Range("L117").Value = 124999.75
$ws.Range("N117").Value = -134177.75
$ws.Range("H122").Value = 1427.2667
$ws.Range("I122").Value = 1314.9286
$ws.Range("K122").Value = 3944.7858
$ws.Range("M122").Value = -1494.7858
$ws.Range("H124").Value = 50429
$ws.Range("J124").Value = 50429
$ws.Range("L124").Value = 50429
$ws.Range("N124").Value = -60249
$ws.Range("H125").Value = 250000
$ws.Range("J125").Value = 250000
$ws.Range("L125").Value = 250000
$ws.Range("N125").Value = -259840
$ws.Range("H136").Value = 2635.6428
$ws.Range("I136").Value = 2635.6428
$ws.Range("K136").Value = 7906.928400000001
$ws.Range("M136").Value = -5356.928400000001
$ws.Range("H139").Value = 71500
$ws.Range("J139").Value = 71500
$ws.Range("L139").Value = 71500
$ws.Range("N139").Value = -81780

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H54").Value = 7133
$ws.Range("I54").Value = 7133
$ws.Range("K54").Value = 7133
$ws.Range("M54").Value = -6649
$ws.Range("H105").Value = 2728.8518
$ws.Range("I105").Value = 2738.3044
$ws.Range("K105").Value = 2738.3044
$ws.Range("M105").Value = -991.3044
$ws.Range("H107").Value = 1086
$ws.Range("I107").Value = 814.7692
$ws.Range("J107").Value = 2849
$ws.Range("K107").Value = 814.7692
$ws.Range("L107").Value = 2849
$ws.Range("M107").Value = 1105.2308
$ws.Range("N107").Value = -6689

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 637.5714
$ws.Range("I16").Value = 637.5714
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 637.5714
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = -350.5714
$ws.Range("N16").ClearContents()
$ws.Range("H58").Value = 68198.92999999999
$ws.Range("I58").Value = 72783.57000000001
$ws.Range("J58").Value = 4014
$ws.Range("K58").Value = 72783.57000000001
$ws.Range("L58").Value = 4014
$ws.Range("M58").Value = -72580.57000000001
$ws.Range("N58").Value = -4420
$ws.Range("H99").Value = 2861.0908
$ws.Range("I99").Value = 1948.5
$ws.Range("J99").Value = 3063.889
$ws.Range("K99").Value = 1948.5
$ws.Range("L99").Value = 3063.889
$ws.Range("M99").Value = -450.5
$ws.Range("N99").Value = -6059.889
$ws.Range("H107").Value = 402.08334
$ws.Range("J107").Value = 452.75
$ws.Range("L107").Value = 452.75
$ws.Range("N107").Value = -4292.75
$ws.Range("H113").Value = 637.5714
$ws.Range("I113").Value = 637.5714
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 637.5714
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = 1532.4286
$ws.Range("N113").ClearContents()
$ws.Range("H126").Value = 2861.0908
$ws.Range("I126").Value = 1948.5
$ws.Range("J126").Value = 3063.889
$ws.Range("K126").Value = 5845.5
$ws.Range("L126").Value = 9191.667000000001
$ws.Range("M126").Value = -3375.5
$ws.Range("N126").Value = -14131.667
$ws.Range("H132").Value = 2388
$ws.Range("I132").Value = 2388
$ws.Range("K132").Value = 7164
$ws.Range("M132").Value = -4634
$ws.Range("H136").Value = 68198.92999999999
$ws.Range("I136").Value = 72783.57000000001
$ws.Range("J136").Value = 4014
$ws.Range("K136").Value = 218350.71
$ws.Range("L136").Value = 12042
$ws.Range("M136").Value = -215800.71
$ws.Range("N136").Value = -17142

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 226.05882
$ws.Range("I12").Value = 156.66667
$ws.Range("J12").Value = 263.9091
$ws.Range("K12").Value = 470.00001
$ws.Range("L12").Value = 791.7273
$ws.Range("M12").Value = -297.00001
$ws.Range("N12").Value = -1137.7273
$ws.Range("H37").Value = 64583.332
$ws.Range("J37").Value = 64583.332
$ws.Range("L37").Value = 193749.996
$ws.Range("N37").Value = -193973.996
$ws.Range("H62").Value = 4514
$ws.Range("J62").Value = 4514
$ws.Range("L62").Value = 13542
$ws.Range("N62").Value = -14914
$ws.Range("H63").Value = 0
$ws.Range("I63").Value = 0
$ws.Range("K63").Value = 0
$ws.Range("M63").ClearContents()
$ws.Range("H65").Value = 4514
$ws.Range("J65").Value = 4514
$ws.Range("L65").Value = 40626
$ws.Range("N65").Value = -47490
$ws.Range("H66").Value = 0
$ws.Range("I66").Value = 0
$ws.Range("K66").Value = 0
$ws.Range("M66").ClearContents()
$ws.Range("H69").Value = 977.9091
$ws.Range("H70").Value = 19999.5
$ws.Range("I70").Value = 19999
$ws.Range("K70").Value = 59997
$ws.Range("M70").Value = -59682
$ws.Range("H72").Value = 977.9091
$ws.Range("H73").Value = 19999.5
$ws.Range("I73").Value = 19999
$ws.Range("K73").Value = 59997
$ws.Range("M73").Value = -58905
$ws.Range("H75").Value = 950
$ws.Range("J75").Value = 950
$ws.Range("L75").Value = 2850
$ws.Range("N75").Value = -4846
$ws.Range("H78").Value = 950
$ws.Range("J78").Value = 950
$ws.Range("L78").Value = 8550
$ws.Range("N78").Value = -18534
$ws.Range("H118").Value = 3008.5
$ws.Range("I118").Value = 3008.5
$ws.Range("K118").Value = 9025.5
$ws.Range("M118").Value = -7782.5
$ws.Range("H123").Value = 51014.5
$ws.Range("J123").Value = 0
$ws.Range("L123").Value = 0
$ws.Range("N123").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 54054.473
$ws.Range("I107").Value = 84449.414
$ws.Range("J107").Value = 1948.8572
$ws.Range("K107").Value = 84449.414
$ws.Range("L107").Value = 1948.8572
$ws.Range("M107").Value = -82529.414
$ws.Range("N107").Value = -5788.8572
$ws.Range("H122").Value = 3084.7778
$ws.Range("I122").Value = 2537.7144
$ws.Range("J122").Value = 4999.5
$ws.Range("K122").Value = 7613.1432
$ws.Range("L122").Value = 14998.5
$ws.Range("M122").Value = -5163.1432
$ws.Range("N122").Value = -19898.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H60").Value = 55000
$ws.Range("J60").Value = 55000
$ws.Range("L60").Value = 55000
$ws.Range("N60").Value = -56018
$ws.Range("H132").Value = 39745.605
$ws.Range("I132").Value = 47189.074
$ws.Range("K132").Value = 141567.222
$ws.Range("M132").Value = -139037.222

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 990.8333
$ws.Range("I100").Value = 913
$ws.Range("J100").Value = 1146.5
$ws.Range("K100").Value = 1826
$ws.Range("L100").Value = 2293
$ws.Range("M100").Value = -1285
$ws.Range("N100").Value = -3375
$ws.Range("H122").Value = 1777.2916
$ws.Range("J122").Value = 1916.6666
$ws.Range("L122").Value = 5749.9998
$ws.Range("N122").Value = -10649.9998
$ws.Range("H132").Value = 34822.88
$ws.Range("I132").Value = 37264.605
$ws.Range("J132").Value = 2470
$ws.Range("K132").Value = 111793.815
$ws.Range("L132").Value = 7410
$ws.Range("M132").Value = -109263.815
$ws.Range("N132").Value = -12470
